$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 150, shifting existing rows 150-220 down to 151-221.
$ws.Rows.Item(150).Insert()

# Populate the newly inserted row 150 with the new data point.
$ws.Cells.Item(150, 1).Value = 8
$ws.Cells.Item(150, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(150, 3).Value = "Coquimbo"
$ws.Cells.Item(150, 4).Value = 45027
$ws.Cells.Item(150, 5).Value = 4
$ws.Cells.Item(150, 6).Value = 100112040
$ws.Cells.Item(150, 7).Value = "Cilantro"
$ws.Cells.Item(150, 8).Value = "Sin especificar"
$ws.Cells.Item(150, 9).Value = "Primera"
$ws.Cells.Item(150, 10).Value = 2200
$ws.Cells.Item(150, 11).Value = 1800
$ws.Cells.Item(150, 12).Value = 2000
$ws.Cells.Item(150, 13).Value = 1900
$ws.Cells.Item(150, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(150, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(150, 16).Value = 1267
$ws.Cells.Item(150, 17).Value = 1.5
$ws.Cells.Item(150, 18).Value = "Hortaliza"
